$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores values as text (e.g. "1.00", "0.420") even
# though they look numeric. Force NumberFormat to Text on each Price cell
# right before writing it so Excel keeps the exact string instead of
# silently coercing it to a number (which would drop trailing zeros).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.900.96"
$ws.Range("E2").Value = "  -0.26%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.524.90"
$ws.Range("E3").Value = "  +0.74%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.28%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "535.66"
$ws.Range("E5").Value = "  -0.37%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.55"
$ws.Range("E6").Value = "  -1.36%  "

# Row 7
$ws.Range("E7").Value = "  -0.33%  "

# Row 8
$ws.Range("E8").Value = "  +1.11%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.525.47"
$ws.Range("E9").Value = "  +0.94%  "

# Row 10
$ws.Range("E10").Value = "  +0.79%  "

# Row 11
$ws.Range("E11").Value = "  -2.68%  "

# Row 12
$ws.Range("E12").Value = "  -1.81%  "

# Row 13
$ws.Range("E13").Value = "  +0.08%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.966.15"
$ws.Range("E14").Value = "  -0.08%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.04"
$ws.Range("E15").Value = "  -0.08%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.954.85"
$ws.Range("E16").Value = "  -0.20%  "

# Row 17
$ws.Range("E17").Value = "  -0.43%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.519.59"
$ws.Range("E18").Value = "  +0.16%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.14"
$ws.Range("E19").Value = "  +1.26%  "

# Row 20
$ws.Range("E20").Value = "  +0.55%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.97"
$ws.Range("E21").Value = "  -0.89%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.06%  "

# Row 23
$ws.Range("E23").Value = "  +2.11%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.22"
$ws.Range("E24").Value = "  +3.49%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.420"
$ws.Range("E25").Value = "  +0.61%  "

# Row 26
$ws.Range("E26").Value = "  -1.46%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.997"
$ws.Range("E27").Value = "  -0.41%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.51"
$ws.Range("E28").Value = "  -2.32%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0770"
$ws.Range("E29").Value = "  -0.12%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.67"
$ws.Range("E30").Value = "  +0.35%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "171.97"
$ws.Range("E31").Value = "  +4.07%  "

# Row 32
$ws.Range("E32").Value = "  -1.32%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.17"
$ws.Range("E33").Value = "  +6.97%  "

# Row 34
$ws.Range("E34").Value = "  -0.19%  "

# Row 35
$ws.Range("E35").Value = "  +2.82%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.41"
$ws.Range("E36").Value = "  -0.13%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.09"
$ws.Range("E37").Value = "  -0.83%  "

# Row 38
$ws.Range("E38").Value = "  -1.40%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.69"
$ws.Range("E39").Value = "  -0.38%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.811"
$ws.Range("E40").Value = "  +1.23%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.59"
$ws.Range("E41").Value = "  -1.19%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "285.15"
$ws.Range("E42").Value = "  +2.21%  "

# Row 43
$ws.Range("E43").Value = "  +0.64%  "

# Row 44
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  -0.20%  "

# Row 45
$ws.Range("E45").Value = "  +2.80%  "

# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "131.43"
$ws.Range("E46").Value = "  +5.05%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.87"
$ws.Range("E47").Value = "  -0.23%  "

# Row 48
$ws.Range("E48").Value = "  -1.40%  "

# Row 49
$ws.Range("E49").Value = "  -0.26%  "

# Row 50
$ws.Range("E50").Value = "  -0.49%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.39"
$ws.Range("E51").Value = "  -1.37%  "
